# Apply the "Changed to USB-C" BOM edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("FLCOS Projector Driver BOM")

# Row 10: resistor designators/value updated (R2,R4,R16 -> R2,R4,R13,R14,R16 ; 20K -> 22K)
$ws.Range("A10").Value = "R2, R4, R13, R14, R16"
$ws.Range("B10").Value = "22K"

# Row 18: manufacturer name for the Schottky diode filled in
$ws.Range("H18").Value = "Panjit"

# Row 29: manufacturer part number text duplicated onto a second line
$ws.Range("I29").Value = "MEM2075-00-140-01-A" + [char]10 + "MEM2075-00-140-01-A"

# Row 30: USB connector swapped from mini-USB (CUI) to USB-C (GCT)
$ws.Range("A30").Value = "USB-C"
$ws.Range("B30").Value = "USB4105-GF-A"
$ws.Range("C30").Value = "USB Type C,2.0"
$ws.Range("D30").Value = "GCT_USB4105-GF-A"
$ws.Range("E30").Value = "USB-C connector"
$ws.Range("F30").Value = "https://www.mouser.com/ProductDetail/GCT/USB4105-GF-A?qs=KUoIvG%2F9IlY%2FMLlBMpStpA%3D%3D"
$ws.Range("H30").Value = "GCT"
$ws.Range("I30").Value = "USB4105-GF-A"
$ws.Range("J30").Value = "640-USB4105-GF-A"

# Row 31: price updated now that the USB-C connector is in the BOM
$ws.Range("K31").Value = 0.81

# Move the active selection to A10, matching the saved workbook's view state
$ws.Range("A10").Select()
